$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.031.62'
$ws.Range("E2").Value = '  -1.39%  '

$ws.Range("D3").Value = '3.564.64'
$ws.Range("E3").Value = '  -1.68%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '''236.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.81%  '

$ws.Range("D6").Value = '''653.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.14%  '

$ws.Range("E7").Value = '  -0.86%  '

$ws.Range("D8").Value = '''0.399'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.42%  '

$ws.Range("E9").Value = '  +0.15%  '

$ws.Range("E10").Value = '  -0.15%  '

$ws.Range("D11").Value = '3.561.84'
$ws.Range("E11").Value = '  -1.61%  '

$ws.Range("E12").Value = '  +0.90%  '

$ws.Range("D13").Value = '''42.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.42%  '

$ws.Range("D14").Value = '''6.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.44%  '

$ws.Range("D15").Value = '4.230.36'
$ws.Range("E15").Value = '  -1.97%  '

$ws.Range("D16").Value = '95.231.95'
$ws.Range("E16").Value = '  -1.07%  '

$ws.Range("D17").Value = '''0.0000253'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.42%  '

$ws.Range("D18").Value = '''8.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.86%  '

$ws.Range("D19").Value = '3.568.49'
$ws.Range("E19").Value = '  -1.69%  '

$ws.Range("D20").Value = '''12.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.39%  '

$ws.Range("D21").Value = '''17.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.06%  '

$ws.Range("D22").Value = '''3.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("D23").Value = '''508.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.54%  '

$ws.Range("D24").Value = '''0.480'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.86%  '

$ws.Range("D25").Value = '''6.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.68%  '

$ws.Range("E26").Value = '  -0.72%  '

$ws.Range("D27").Value = '''95.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.37%  '

$ws.Range("D28").Value = '''12.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.66%  '

$ws.Range("D29").Value = '3.759.69'
$ws.Range("E29").Value = '  -1.45%  '

$ws.Range("E30").Value = '  -3.65%  '

$ws.Range("D31").Value = '''0.144'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.93%  '

$ws.Range("D32").Value = '''11.55'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.26%  '

$ws.Range("D33").Value = '''0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("D34").Value = '''0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("B35").Value = 'Cronos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D35").Value = '''0.176'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.15%  '

$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '''31.87'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.82%  '

$ws.Range("D37").Value = '''1.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +16.43%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").Value = '''8.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.35%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '''603.38'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.99%  '

$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = '''0.560'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.64%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '''0.151'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.48%  '

$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").Value = '''1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("D43").Value = '''0.905'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.41%  '

$ws.Range("E44").Value = '  +6.10%  '

$ws.Range("D45").Value = '''34.81'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +29.29%  '

$ws.Range("D46").Value = '''5.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.68%  '

$ws.Range("E47").Value = '  +3.44%  '

$ws.Range("D48").Value = '''23.42'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.62%  '

$ws.Range("D49").Value = '''0.0416'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.23%  '

$ws.Range("D50").Value = '''3.50'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.61%  '

$ws.Range("D51").Value = '''8.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.29%  '

